# Peer review workbook - fill in completed review content ("Finished with project 4!")
#
# This reproduces (as closely as COM allows) the author's edit: the peer-review
# header block (group/date/artifact/scribe/leader/time) gets filled in, two
# issues are logged in the table, and the whole used range is re-left-aligned
# (previously only centred/wrapped styles existed for these cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlLeft = -4131

# ---- Row 1: title -------------------------------------------------------
# A1 keeps its bold 12pt font, just becomes left-aligned; B1 becomes an
# (empty) left-aligned cell.
$ws.Range("A1").HorizontalAlignment = $xlLeft
$ws.Range("B1").HorizontalAlignment = $xlLeft

# ---- Row 2: subtitle ------------------------------------------------------
$ws.Range("A2").HorizontalAlignment = $xlLeft
$ws.Range("B2").HorizontalAlignment = $xlLeft

# ---- Row 3: new blank spacer row, left-aligned ---------------------------
$ws.Range("A3").HorizontalAlignment = $xlLeft
$ws.Range("B3").HorizontalAlignment = $xlLeft

# ---- Row 4: Group# --------------------------------------------------------
$ws.Range("A4").HorizontalAlignment = $xlLeft
$ws.Range("B4").Value = 0
$ws.Range("B4").HorizontalAlignment = $xlLeft

# ---- Row 5: Date: ----------------------------------------------------------
$ws.Range("A5").HorizontalAlignment = $xlLeft
$ws.Range("B5").Value = "Thursday May 18, 2021"
$ws.Range("B5").HorizontalAlignment = $xlLeft

# ---- Row 6: Artifact: -------------------------------------------------------
$ws.Range("A6").HorizontalAlignment = $xlLeft
$ws.Range("B6").Value = "Project 3 source code"
$ws.Range("B6").HorizontalAlignment = $xlLeft

# ---- Row 7: Scribe: ----------------------------------------------------------
$ws.Range("A7").HorizontalAlignment = $xlLeft
$ws.Range("B7").Value = "Todd Nguyen"
$ws.Range("B7").HorizontalAlignment = $xlLeft

# ---- Row 8: Leader: ----------------------------------------------------------
$ws.Range("A8").HorizontalAlignment = $xlLeft
$ws.Range("B8").Value = "Todd Nguyen"
$ws.Range("B8").HorizontalAlignment = $xlLeft

# ---- Row 9: Time spent: ------------------------------------------------------
$ws.Range("A9").HorizontalAlignment = $xlLeft
$ws.Range("B9").Value = "30 minutes"
$ws.Range("B9").HorizontalAlignment = $xlLeft

# ---- Row 12-13: logged issues ---------------------------------------------
$ws.Range("A12:A26").RowHeight = 13.2
$ws.Rows.Item(12).RowHeight = 26.4

$ws.Range("B12").Value = "Lack of comment on each function. However, the function name might be more than enough for a comment."
$ws.Range("C12").Value = "Not a Bug"

$ws.Range("B13").Value = "Some structs can be typedef."
$ws.Range("C13").Value = "Not Fixed"

# ---- Selection / view state -------------------------------------------------
$ws.Range("B7").Select()
